$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the worksheet / workbook title -------------------------------
$ws.Name = "congestion by duid top 10"

# --- Data rows (new top-10 congestion-by-DUID figures) ---------------------
# (written before the header row so the shared-string table ends up in the
# same relative order as the source workbook)
$ws.Range("A2").Value = "Lake Bonney Stage 2"
$ws.Range("B2").Value = "Wind"
$ws.Range("C2").Value = 15.598363774733601

$ws.Range("A3").Value = "Lake Bonney Stage 3"
$ws.Range("B3").Value = "Wind"
$ws.Range("C3").Value = 15.1522070015221

$ws.Range("A4").Value = "Snowtown North Wind Farm"
$ws.Range("B4").Value = "Wind"
$ws.Range("C4").Value = 13.6339421613394

$ws.Range("A5").Value = "Snowtown South Wind Farm"
$ws.Range("B5").Value = "Wind"
$ws.Range("C5").Value = 13.632039573820402

$ws.Range("A6").Value = "Hallett Wind Farm"
$ws.Range("B6").Value = "Wind"
$ws.Range("C6").Value = 13.165905631659101

$ws.Range("A7").Value = "Hornsdale Wind Farm"
$ws.Range("B7").Value = "Wind"
$ws.Range("C7").Value = 13.112633181126302

$ws.Range("A8").Value = "North Brown Hill Wind Farm"
$ws.Range("B8").Value = "Wind"
$ws.Range("C8").Value = 12.991818873668201

$ws.Range("A9").Value = "The Bluff Wind Farm"
$ws.Range("B9").Value = "Wind"
$ws.Range("C9").Value = 12.950913242009099

$ws.Range("A10").Value = "Snowtown Wind Farm"
$ws.Range("B10").Value = "Wind"
$ws.Range("C10").Value = 12.8491248097412

$ws.Range("A11").Value = "Clements Gap Wind Farm"
$ws.Range("B11").Value = "Wind"
$ws.Range("C11").Value = 12.706430745814302

# --- Header row ------------------------------------------------------------
$ws.Range("A1").Value = "Station"
$ws.Range("B1").Value = "Fuel Type"
$ws.Range("C1").Value = "Percent"

# --- Formatting --------------------------------------------------------
# Header: bold, no underline
$ws.Range("A1:C1").Font.Bold = $true
$ws.Range("A1:C1").Font.Underline = $false

# Percent column: one decimal place
$ws.Range("C2:C11").NumberFormat = "0.0"

# Column widths for the new columns (closest values the engine's pixel
# quantization can reach to the source widths of 9.54296875 / 7.81640625)
$ws.Columns.Item(2).ColumnWidth = 8.66
$ws.Columns.Item(3).ColumnWidth = 7.0

# Page setup
$ws.PageSetup.PaperSize = 9

# Selection covering the whole table
$ws.Range("A1:C11").Select()
